# Apply the edits described by the commit:
#  - Update the two "next" time-slot strings on Sheet1 (B8/B9)
#    from 19:45-19:49 / 19:50-19:54 to 20:05-20:09 / 20:10-20:14
#  - Move the active selection on Sheet1 from A13 to A12
#  - Nudge the workbook window position/size to match the new view
#    (best-effort; harmless if the host does not persist it)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the shared-string backed cell values.
$ws.Range("B8").Value = "20:05 - 20:09"
$ws.Range("B9").Value = "20:10 - 20:14"

# Move the selection/active cell from A13 to A12.
$ws.Range("A12").Select()

# Match the updated workbook window geometry from the diff.
$excel.ActiveWindow.Left = 2145
$excel.ActiveWindow.Top = 465
$excel.ActiveWindow.Width = 24360
$excel.ActiveWindow.Height = 15120
